$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark IA2 tasks (C16, C17, H20, H21) as "Completed" instead of
# "Incomplete". Re-use the existing "Completed"/Good formatting already
# present on other cells (e.g. C4) by copying its format over, so the
# workbook's shared style index is reused rather than a new one being
# created.
$targets = @("C16", "C17", "H20", "H21")
foreach ($addr in $targets) {
    $ws.Range("C4").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range($addr).Value = "Completed"
}
$excel.CutCopyMode = 0

# Update the selected cell to match the author's last selection
$ws.Range("C17").Select() | Out-Null
